$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, shifting existing rows 6-115 down to 7-116
$ws.Rows.Item(6).Insert()

# Populate the new row 6 with its data
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Macroferia Regional de Talca"
$ws.Range("C6").Value = "Maule"
$ws.Range("D6").Value = 44882
$ws.Range("E6").Value = 7
$ws.Range("F6").Value = 100112022
$ws.Range("G6").Value = "Arveja Verde"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 500
$ws.Range("K6").Value = 17000
$ws.Range("L6").Value = 17000
$ws.Range("M6").Value = 17000
$ws.Range("N6").Value = "`$/saco 25 kilos"
$ws.Range("O6").Value = "Región del Maule"
$ws.Range("P6").Value = 680
$ws.Range("Q6").Value = 25
$ws.Range("R6").Value = "Hortaliza"
